# Updated cryptos list on Thu Nov 28 17:37:56 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the 50 coin rows (2-51) on the active sheet, matching the latest pull
# from coinranking.com.
#
# Several Price cells hold digit strings that Excel's Range.Value setter
# would otherwise auto-interpret as numbers (e.g. dropping the trailing
# zero in "6.50" -> 6.5, or silently changing the cell's stored type from
# text to a number even when the digits round-trip, e.g. "650.04"). Those
# are entered with a leading apostrophe so they're kept as literal text,
# exactly like typing them into Excel by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.166.76'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").Value = '3.561.56'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  -1.31%  '
$ws.Range("D6").Value = '''650.04'
$ws.Range("E6").Value = '  +1.93%  '
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("D11").Value = '3.561.35'
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("D13").Value = '''42.32'
$ws.Range("E13").Value = '  -2.50%  '
$ws.Range("D14").Value = '''6.50'
$ws.Range("E14").Value = '  +0.72%  '
$ws.Range("D15").Value = '4.220.58'
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("D16").Value = '95.048.41'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").Value = '3.571.74'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("E19").Value = '  -2.33%  '
$ws.Range("D20").Value = '''12.64'
$ws.Range("E20").Value = '  -3.36%  '
$ws.Range("D21").Value = '''17.90'
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("E22").Value = '  +1.59%  '
$ws.Range("D23").Value = '''505.82'
$ws.Range("E23").Value = '  -2.15%  '
$ws.Range("D24").Value = '''0.477'
$ws.Range("E24").Value = '  -5.81%  '
$ws.Range("D25").Value = '''6.75'
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = '''95.21'
$ws.Range("E27").Value = '  -1.88%  '
$ws.Range("D28").Value = '''12.46'
$ws.Range("E28").Value = '  +0.89%  '
$ws.Range("D29").Value = '3.752.41'
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").Value = '''3.02'
$ws.Range("E30").Value = '  -2.91%  '
$ws.Range("D31").Value = '''11.41'
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").Value = '''0.142'
$ws.Range("E32").Value = '  -3.35%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").Value = '''0.999'
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("E35").Value = '  -3.66%  '
$ws.Range("D36").Value = '''31.74'
$ws.Range("E36").Value = '  +5.04%  '
$ws.Range("D37").Value = '''0.558'
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("D38").Value = '''8.41'
$ws.Range("E38").Value = '  +6.94%  '
$ws.Range("E39").Value = '  +6.93%  '
$ws.Range("D40").Value = '''576.64'
$ws.Range("E40").Value = '  -0.43%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").Value = '''0.150'
$ws.Range("E42").Value = '  -0.67%  '
$ws.Range("D43").Value = '''0.901'
$ws.Range("E43").Value = '  -2.31%  '
$ws.Range("D44").Value = '''1.76'
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("D45").Value = '''2.28'
$ws.Range("E45").Value = '  +4.54%  '
$ws.Range("D46").Value = '''5.68'
$ws.Range("E46").Value = '  +0.92%  '
$ws.Range("E47").Value = '  -1.92%  '
$ws.Range("D48").Value = '''0.0412'
$ws.Range("E48").Value = '  -4.76%  '
$ws.Range("D49").Value = '''3.59'
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("D50").Value = '''33.18'
$ws.Range("E50").Value = '  +31.32%  '
$ws.Range("D51").Value = '''53.20'
$ws.Range("E51").Value = '  -1.31%  '
